# Update technology catalog values and refresh the "current" view state
# (active sheet/selection) to match the latest working state - "A4 for jonathan".

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("inputdisp")
$wsPipe  = $wb.Worksheets.Item("endofpipe")

# --- inputdisp: t1 (row 2) / t2 (row 3) values ---
$wsInput.Range("B2").Value = 2.5
$wsInput.Range("C2").Value = 0.6
$wsInput.Range("E2").Value = 0.35
$wsInput.Range("F2").Value = 0.01
$wsInput.Range("G2").Value = 3.5
$wsInput.Range("H2").Value = 0.98

$wsInput.Range("B3").Value = 5
$wsInput.Range("E3").Value = 0.35
$wsInput.Range("F3").Value = 0.001
$wsInput.Range("G3").Value = 7
$wsInput.Range("H3").Value = 0.98

# --- endofpipe: t1 (row 2) / t2 (row 3) values ---
$wsPipe.Range("B2").Value = 3
$wsPipe.Range("E2").Value = 0.1
$wsPipe.Range("F2").Value = 0.00001
$wsPipe.Range("G2").Value = 5

$wsPipe.Range("B3").Value = 10
$wsPipe.Range("E3").Value = 0.7
$wsPipe.Range("F3").Value = 0.000002
$wsPipe.Range("G3").Value = 13.5

# --- View state: leave a selection on inputdisp, then make endofpipe the
#     active (selected) sheet/tab with its own selection ---
$wsInput.Range("B31").Select()

$wsPipe.Activate()
$wsPipe.Range("D22").Select()

$wb.Save()
